$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-27 from 45208 (2023-10-09)
# to 45212 (2023-10-13), keeping existing cell formatting.
for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45208) {
        $cell.Value = 45212
    }
}
